# 14 Switch Statement -> 15 Calculator (chapter renumbering + topic rename)
# and hyperlink index bump (index=14 -> index=15, plus new video id),
# plus a red outline added to the second picture on slide 6.

$p = $ppt.ActivePresentation

# --- Slide 1 (title slide) ---------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "15 Calculator"

# --- Slide 2 (section title slide) --------------------------------------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "15 Calculator"
$s2.Shapes.Item(2).TextFrame.TextRange.Text = "Calculator"
$s2.Shapes.Item(3).TextFrame.TextRange.Text = "https://www.youtube.com/watch?v=d6dnCQS8DCk&list=PL0eyrZgxdwhwBToawjm9faF1ixePexft-&index=15"

# --- Slide 3 (14.1 Code -> 15.1 Code) -----------------------------------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "15.1 Code"

# --- Slide 4 (14.1 Code page, title split across two runs: "14" + ".1 Code") ---
$s4 = $p.Slides.Item(4)
$titleRange4 = $s4.Shapes.Item(1).TextFrame.TextRange
$titleRange4.Characters(1, 2).Text = "15"
$s4.Shapes.Item(3).TextFrame.TextRange.Text = "https://www.youtube.com/watch?v=d6dnCQS8DCk&list=PL0eyrZgxdwhwBToawjm9faF1ixePexft-&index=15"

# --- Slide 5 (14.2 Verify -> 15.2 Verify) -------------------------------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "15.2 Verify"

# --- Slide 6 (14.2 Verify page) -----------------------------------------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "15.2 Verify"
$s6.Shapes.Item(3).TextFrame.TextRange.Text = "https://www.youtube.com/watch?v=d6dnCQS8DCk&list=PL0eyrZgxdwhwBToawjm9faF1ixePexft-&index=15"

# Add a red (C00000) outline to the second picture ("Picture 6") on slide 6
$pic6 = $s6.Shapes.Item(7)
$pic6.Line.ForeColor.RGB = 0x0000C0
